$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 41436
$ws.Range("A11").NumberFormat = 'ddd\ dd/mm/yyyy'

$ws.Range("B11").Value = 1.5

$ws.Range("D11").Value = "Implementation of tc11_mutex"

$ws.Range("C11").Select()
